$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts old D:K to E:L)
$ws.Columns("D:D").Insert()

# Copy formatting from the (now shifted) old-D column (now column E) into new column D
# so the new column matches the existing per-row number formats/styles.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the brand-new column D with the latest quarter's figures
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 7100
$ws.Range("D9").Value = 6000
$ws.Range("D10").Value = 1100
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 10700
$ws.Range("D18").Value = -3600
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = -3400
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = -3500
$ws.Range("D24").Value = 200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -3700
$ws.Range("D27").Value = -3700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = -3700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -3700
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 1000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 3000
$ws.Range("D44").Value = 3700
$ws.Range("D45").Value = 2800
$ws.Range("D46").Value = 10500
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 1100
$ws.Range("D49").Value = 16500
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 600
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 28700
$ws.Range("D57").Value = 4000
$ws.Range("D58").Value = 3500
$ws.Range("D59").Value = 4600
$ws.Range("D60").Value = 12100
$ws.Range("D61").Value = 3000
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 15100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -19200
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 13600
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = -3700
$ws.Range("D83").Value = 100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 600
$ws.Range("D91").Value = -100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -1100
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 1100
$ws.Range("D101").Value = -100
$ws.Range("D102").Value = 600

# Two cells that also needed a data correction to "NA" at the same time as the insert
$ws.Range("I91").Value = "NA"
$ws.Range("J91").Value = "NA"
